# "Finished up reading from week 5"
#
# - Fill in the "Actual time length to complete" values (column C) for the
#   tasks in week5 that were still blank (rows 4-8), entered as hours:minutes
#   (Excel stores these as a fraction of a day).
# - Move the active window/selection from week6 (last thing looked at, cell
#   D15) to week5 (cell C9), since that is the sheet now being worked on.

$wb = $excel.ActiveWorkbook

# Reposition the application window (best effort - matches the author moving
# the Excel window across displays between edits).
$win = $excel.ActiveWindow
$win.Left = -27580
$win.Top = 500

$ws5 = $wb.Worksheets.Item("week5")

# Actual time length to complete, entered as h:mm and stored by Excel as
# minutes / 60 / 24 of a day.
$ws5.Range("C4").Value = 23 / 60 / 24   # 0:23
$ws5.Range("C5").Value = 45 / 60 / 24   # 0:45
$ws5.Range("C6").Value = 60 / 60 / 24   # 1:00
$ws5.Range("C7").Value = 10 / 60 / 24   # 0:10
$ws5.Range("C8").Value = 1 / 60 / 24    # 0:01

# Make week5 the active sheet/tab with C9 selected - this also clears the
# "tabSelected"/selection state that used to be on week6.
$ws5.Activate()
$ws5.Range("C9").Select()
